$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 3 and row 4 in columns B, C, D
$tmpB = $ws.Range("B3").Value2
$tmpC = $ws.Range("C3").Value2
$tmpD = $ws.Range("D3").Value2

$ws.Range("B3").Value2 = $ws.Range("B4").Value2
$ws.Range("C3").Value2 = $ws.Range("C4").Value2
$ws.Range("D3").Value2 = $ws.Range("D4").Value2

$ws.Range("B4").Value2 = $tmpB
$ws.Range("C4").Value2 = $tmpC
$ws.Range("D4").Value2 = $tmpD

# Update the selected cell
$ws.Range("A5").Select()
